$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.618.20"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "2.286.20"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'95.45"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").Value = "'268.04"
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("E7").Value = "  -1.12%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -2.68%  "
$ws.Range("D10").Value = "'45.28"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").Value = "'0.0931"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "'7.88"
$ws.Range("E12").Value = "  -3.60%  "
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "2.627.31"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").Value = "'15.30"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "'0.847"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "2.286.85"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").Value = "43.493.74"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").Value = "'0.0000108"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("D21").Value = "'71.94"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "'2.57"
$ws.Range("E22").Value = "  +11.91%  "
$ws.Range("D23").Value = "'232.19"
$ws.Range("E23").Value = "  -2.86%  "
$ws.Range("D24").Value = "'9.14"
$ws.Range("E24").Value = "  -5.46%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.59"
$ws.Range("E25").Value = "  +2.81%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "'11.19"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  +2.40%  "
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("E30").Value = "  -6.16%  "
$ws.Range("D31").Value = "'174.85"
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("D32").Value = "'21.74"
$ws.Range("E32").Value = "  -3.79%  "
$ws.Range("D33").Value = "'0.0891"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("E34").Value = "  -4.04%  "
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("E36").Value = "  -3.68%  "
$ws.Range("D37").Value = "'0.0350"
$ws.Range("E37").Value = "  -2.99%  "
$ws.Range("E38").Value = "  -3.22%  "
$ws.Range("E39").Value = "  -4.08%  "
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").Value = "'2.31"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'12.29"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'65.44"
$ws.Range("E43").Value = "  +7.83%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("E46").Value = "  -0.55%  "
$ws.Range("E47").Value = "  -6.13%  "
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").Value = "'96.33"
$ws.Range("E49").Value = "  -4.12%  "
$ws.Range("D50").Value = "'0.432"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.507.08"
$ws.Range("E51").Value = "  -1.18%  "
